$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Schedules")

# --- Section title: "20-Week variant" ---
$ws.Range("A36").Value = "20-Week variant"
$ws.Range("A36").Font.Bold = $true
$ws.Range("A36").Font.Size = 14

# --- Header row for the new table ---
$ws.Range("A37").Value = "Week"
$ws.Range("B37").Value = "Chapters"
$ws.Range("C37").Value = "Subjects"
$ws.Range("D37").Value = "Remarks"
$ws.Range("A37:D37").Font.Bold = $true
$ws.Range("A37:D37").Font.Size = 12

# --- Week numbers (column A) for all 20 data rows ---
$ws.Cells.Item(38, 1).Value = 1
$ws.Cells.Item(39, 1).Value = 2
$ws.Cells.Item(40, 1).Value = 3
$ws.Cells.Item(41, 1).Value = 4
$ws.Cells.Item(42, 1).Value = 5
$ws.Cells.Item(43, 1).Value = 6
$ws.Cells.Item(44, 1).Value = 7
$ws.Cells.Item(45, 1).Value = 8
$ws.Cells.Item(46, 1).Value = 9
$ws.Cells.Item(47, 1).Value = 10
$ws.Cells.Item(48, 1).Value = 11
$ws.Cells.Item(49, 1).Value = 12
$ws.Cells.Item(50, 1).Value = 13
$ws.Cells.Item(51, 1).Value = 14
$ws.Cells.Item(52, 1).Value = 15
$ws.Cells.Item(53, 1).Value = 16
$ws.Cells.Item(54, 1).Value = 17
$ws.Cells.Item(55, 1).Value = 18
$ws.Cells.Item(56, 1).Value = 19
$ws.Cells.Item(57, 1).Value = 20

# --- Weeks 2-4: brand-new subject matter (CH1 - CH4) ---
$ws.Range("B39").Value = "CH1"
$ws.Range("C39").Value = "Basics of Java and JVM"

$ws.Range("B40").Value = "CH2 & CH3"
$ws.Range("C40").Value = "Working with variables & values and control flow"

$ws.Range("B41").Value = "CH4"
$ws.Range("C41").Value = "Java Core APIs (String, Date, Math)"

# --- Weeks 5-13: re-use the same chapters/subjects as the 10-week variant ---
$ws.Range("B42").Value = "Review CH1-4"
$ws.Range("C42").Value = "Classes, Inheritance, Overriding and Hiding"

$ws.Range("B43").Value = "CH5"
$ws.Range("C43").Value = "Methods, Static/Nonstatic, Overloading"

$ws.Range("B44").Value = "CH6"
$ws.Range("C44").Value = "Classes, Inheritance, Overriding and Hiding"

$ws.Range("B45").Value = "CH7"
$ws.Range("C45").Value = "Interfaces, Enums, Records, Nesting"

$ws.Range("B46").Value = "CH8"
$ws.Range("C46").Value = "Lambdas"

$ws.Range("B47").Value = "Review CH5-8"

$ws.Range("B48").Value = "CH9"
$ws.Range("C48").Value = "Collections and Generics"

$ws.Range("B49").Value = "CH10"
$ws.Range("C49").Value = "Streams (and streaming operations)"

$ws.Range("B50").Value = "Review CH9-10"

# --- Weeks 14-15: brand-new subject matter (CH11 - CH12) ---
$ws.Range("B51").Value = "CH11"
$ws.Range("C51").Value = "Exceptions, Localization, Resource Management"

$ws.Range("B52").Value = "CH12"
$ws.Range("C52").Value = "Modules"

# --- Week 20: final review ---
$ws.Range("B57").Value = "Review CH1-15"

# --- Week 1: kick-off / planning session ---
$ws.Range("B38").Value = "-"
$ws.Range("C38").Value = "Kick-off / Planning"
$ws.Range("D38").Value = "Check if everyone has received a book and Enthuware licenses. If necessary, assist with installing Enthuware. Make planning. Shorter session (<2hr)"

# --- Week 16: review ---
$ws.Range("B53").Value = "Review CH11-12"

# --- Weeks 17-18: re-use the same chapters/subjects/remarks as the 10-week variant ---
$ws.Range("B54").Value = "CH13"
$ws.Range("C54").Value = "Concurrency"
$ws.Range("D54").Value = "Concurrency usually needs extra attention, because the mental models are harder. Make sure to use plenty of examples."

$ws.Range("B55").Value = "CH14"
$ws.Range("C55").Value = "I/O, Files, Database access"

# --- Week 19: review ---
$ws.Range("B56").Value = "Review CH13-14"

# --- Selection mirrors the freshly-added block ---
$ws.Range("A36:D57").Select()
